$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (losing a significant trailing zero) -- force Text format first, then
# restore the default "Normal" style so no stray formatting is left behind.
$textForceCells = @("D27", "D30", "D35", "D42", "D50")
foreach ($cell in $textForceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# --- Price (column D) updates ---
$ws.Range("D2").Value = "22.354.26"
$ws.Range("D3").Value = "1.562.91"
$ws.Range("D6").Value = "289.45"
$ws.Range("D7").Value = "0.3699"
$ws.Range("D8").Value = "49.29"
$ws.Range("D9").Value = "0.3377"
$ws.Range("D10").Value = "1.163"
$ws.Range("D11").Value = "0.07635"
$ws.Range("D13").Value = "21.37"
$ws.Range("D14").Value = "6.038"
$ws.Range("D15").Value = "6.916"
$ws.Range("D16").Value = "1.558.37"
$ws.Range("D18").Value = "90.05"
$ws.Range("D19").Value = "0.06725"
$ws.Range("D21").Value = "6.243"
$ws.Range("D23").Value = "0.5294"
$ws.Range("D25").Value = "22.355.14"
$ws.Range("D26").Value = "2.376"
$ws.Range("D27").Value = "2.790"
$ws.Range("D28").Value = "20.13"
$ws.Range("D29").Value = "145.71"
$ws.Range("D30").Value = "4.980"
$ws.Range("D31").Value = "125.44"
$ws.Range("D32").Value = "1.737.33"
$ws.Range("D33").Value = "6.202"
$ws.Range("D34").Value = "2.006"
$ws.Range("D35").Value = "1.000"
$ws.Range("D37").Value = "0.08435"
$ws.Range("D39").Value = "0.2319"
$ws.Range("D40").Value = "5.528"
$ws.Range("D41").Value = "0.06425"
$ws.Range("D42").Value = "1.290"
$ws.Range("D44").Value = "0.6326"
$ws.Range("D45").Value = "14.16"
$ws.Range("D48").Value = "3.755"
$ws.Range("D49").Value = "2.096"
$ws.Range("D50").Value = "1.260"
$ws.Range("D51").Value = "124.16"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -4.52%  "
$ws.Range("E3").Value = "  -5.05%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("E7").Value = "  -2.59%  "
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("E9").Value = "  -3.28%  "
$ws.Range("E10").Value = "  -4.56%  "
$ws.Range("E11").Value = "  -5.38%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  -3.07%  "
$ws.Range("E14").Value = "  -4.39%  "
$ws.Range("E15").Value = "  -4.98%  "
$ws.Range("E16").Value = "  -5.55%  "
$ws.Range("E17").Value = "  -7.23%  "
$ws.Range("E18").Value = "  -5.10%  "
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  -5.87%  "
$ws.Range("E22").Value = "  -4.88%  "
$ws.Range("E23").Value = "  -7.62%  "
$ws.Range("E24").Value = "  -3.27%  "
$ws.Range("E25").Value = "  -4.58%  "
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("E27").Value = "  -6.09%  "
$ws.Range("E28").Value = "  -4.32%  "
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("E33").Value = "  -9.27%  "
$ws.Range("E34").Value = "  -5.93%  "
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("E36").Value = "  -10.79%  "
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("E40").Value = "  -6.05%  "
$ws.Range("E41").Value = "  -6.04%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  -8.84%  "
$ws.Range("E44").Value = "  -7.41%  "
$ws.Range("E45").Value = "  -8.84%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  -6.04%  "
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("E49").Value = "  -6.51%  "
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("E51").Value = "  -2.43%  "

# Restore default styling on the cells we temporarily forced to Text format
foreach ($cell in $textForceCells) {
    $ws.Range($cell).Style = "Normal"
}

Write-Output "cryptos list updated"
